# The edit swaps the content of rows 63/64 (Craterellus lutescens <->
# Sarcodon imbricatus records) and rows 65/66 (Lactarius scrobiculatus <->
# Sanicula europaea records). Only the cells that actually differ between
# the two rows are touched, so unrelated/unchanged cells (dates, K/N/T-W,
# AD/AE/AG, AW/AX, ...) are left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 63 <-> Row 64 -------------------------------------------------
$ws.Range("A63").Value = 111683853
$ws.Range("B63").Value = 90687
$ws.Range("E63").Value = 5964
$ws.Range("F63").Value = "Fjällig taggsvamp s.str."
$ws.Range("G63").Value = "Sarcodon imbricatus s.str."
$ws.Range("H63").Value = "(L.:Fr.) P.Karst."
$ws.Range("Z63").Value = "09:34"
$ws.Range("AB63").Value = "09:34"

$ws.Range("A64").Value = 111683845
$ws.Range("B64").Value = 89183
$ws.Range("E64").Value = 3215
$ws.Range("F64").Value = "Rödgul trumpetsvamp"
$ws.Range("G64").Value = "Craterellus lutescens"
$ws.Range("H64").Value = "(Fr.) Fr."
$ws.Range("Z64").Value = "09:36"
$ws.Range("AB64").Value = "09:36"

# --- Row 65 <-> Row 66 -------------------------------------------------
$ws.Range("A65").Value = 111683856
$ws.Range("B65").Value = 108219
$ws.Range("E65").Value = 219711
$ws.Range("F65").Value = "Sårläka"
$ws.Range("G65").Value = "Sanicula europaea"
$ws.Range("H65").Value = "L."
$ws.Range("I65").ClearContents()
$ws.Range("J65").ClearContents()
$ws.Range("P65").Value = "Fiskarsundet, Srm"
$ws.Range("Q65").Value = 689111.5690902721
$ws.Range("R65").Value = 6570305.953062683
$ws.Range("S65").Value = 23
$ws.Range("Z65").Value = "09:34"
$ws.Range("AB65").Value = "09:34"

$ws.Range("A66").Value = 111683850
$ws.Range("B66").Value = 90332
$ws.Range("E66").Value = 4769
$ws.Range("F66").Value = "Svavelriska"
$ws.Range("G66").Value = "Lactarius scrobiculatus"
$ws.Range("H66").Value = "(Scop.:Fr.) Fr."
$ws.Range("I66").NumberFormat = "@"
$ws.Range("I66").Value = "3"
$ws.Range("J66").Value = "fruktkroppar"
$ws.Range("P66").Value = "Bergaholm, Tyresö kn, Srm"
$ws.Range("Q66").Value = 689075.4602011892
$ws.Range("R66").Value = 6570319.534944151
$ws.Range("S66").Value = 20
$ws.Range("Z66").Value = "09:25"
$ws.Range("AB66").Value = "09:25"
